# CourseIntro.pptx — small text updates on the "Resources" slide (slide 15).
#
# The "Content Placeholder 4" shape holds a 3x2 resources table:
#   row2,col2: "Student Hours*"               -> "Drop-In Hours*"
#   row3,col1: "Email Professor and/or your TA" -> "Email Cameron and/or Jiaxin"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)

$shp = $s.Shapes("Content Placeholder 4")
$tbl = $shp.Table

$tbl.Cell(2, 2).Shape.TextFrame.TextRange.Text = "Drop-In Hours*"
$tbl.Cell(3, 1).Shape.TextFrame.TextRange.Text = "Email Cameron and/or Jiaxin"
